# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value for "展览" sheet
$exhibitionUpdates = @{
    3  = 21814
    5  = 8149
    7  = 331
    10 = 216
    11 = 552
    13 = 251
    14 = 1034
    15 = 1376
    17 = 93
    22 = 378
    23 = 1244
    24 = 86
    27 = 5212
    30 = 176
    31 = 5297
    32 = 40
    34 = 78
    36 = 13589
    38 = 169
    41 = 358
    42 = 499
    43 = 4112
    44 = 57
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F-column value for "全部类型" sheet
$allUpdates = @{
    3  = 21814
    5  = 8149
    7  = 331
    10 = 216
    11 = 552
    13 = 251
    14 = 1034
    15 = 1376
    17 = 93
    22 = 378
    23 = 1244
    24 = 86
    28 = 5212
    32 = 176
    34 = 5297
    35 = 40
    37 = 78
    39 = 13589
    41 = 169
    44 = 358
    45 = 499
    46 = 4112
    47 = 57
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
